$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '68.962.79'
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  -0.06%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.929.86'
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  +4.89%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.04%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '604.79'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.67%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '165.38'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -0.38%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '3.927.04'
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +4.88%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.05%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  -1.33%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.167'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -3.75%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '6.41'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -0.04%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +0.70%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '37.28'
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -1.30%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  -0.92%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '4.588.78'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +4.91%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '3.922.69'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +4.53%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '69.085.85'
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -0.04%  '

# Row 18
$ws.Cells.Item(18, 5).Value = '  +0.64%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  -0.45%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '17.16'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -2.68%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '11.15'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -1.05%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '489.02'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -0.41%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.727'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +0.19%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  +12.36%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '84.54'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -0.05%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  -0.44%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '12.14'
$ws.Cells.Item(27, 4).Style = "Normal"

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '10.12'
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +0.71%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +0.07%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  -0.50%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '4.081.65'
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +4.95%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'ImmutableX'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '2.40'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -1.71%  '

# Row 33
$ws.Cells.Item(33, 2).Value = 'NEARProtocol'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '7.88'
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -3.59%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '32.48'
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +2.45%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '3.877.51'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +5.29%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  -0.39%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  +2.92%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  +1.85%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '5.94'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -0.04%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +0.11%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -0.90%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'Bittensor'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '441.95'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +3.41%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'dogwifhat'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '3.02'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -3.75%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  +0.23%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '48.49'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -0.28%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '8.48'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +0.36%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '26.88'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +15.06%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '2.854.11'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +2.18%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '142.09'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +0.51%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  +1.88%  '
